$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was updated
# for every data row (2 through 176) from 45224 (2023-10-25) to
# 45233 (2023-11-03). Update the whole range in one shot using the
# Excel date serial number so the existing date formatting (style)
# on these cells is left untouched.
$ws.Range("C2:C176").Value = 45233
